# Apply updates to the "Fruta, Macroferia Regional de Talca - Granada" sheet.
# Each row 2-21 in the worksheet gets specific cell values corrected
# (dates, quality, volume, prices, units, origin, price/kg, kg/unit)
# per the authoritative dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = [datetime]"2021-05-12"
$ws.Range("M2").Value = 250
$ws.Range("D3").Value = [datetime]"2021-06-11"
$ws.Range("L3").Value = 'Especial'
$ws.Range("M3").Value = 150
$ws.Range("N3").Value = 18000
$ws.Range("O3").Value = 18000
$ws.Range("P3").Value = 18000
$ws.Range("S3").Value = 1000
$ws.Range("D4").Value = [datetime]"2021-06-11"
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 100
$ws.Range("D5").Value = [datetime]"2021-04-30"
$ws.Range("L5").Value = 'Especial'
$ws.Range("M5").Value = 300
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 20000
$ws.Range("R5").Value = 'Provincia de Limarí'
$ws.Range("S5").Value = 1111
$ws.Range("D6").Value = [datetime]"2021-04-08"
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 50
$ws.Range("N6").Value = 12000
$ws.Range("O6").Value = 12000
$ws.Range("P6").Value = 12000
$ws.Range("Q6").Value = '$/caja 15 kilos granel'
$ws.Range("R6").Value = 'Región Metropolitana'
$ws.Range("S6").Value = 800
$ws.Range("T6").Value = 15
$ws.Range("D7").Value = [datetime]"2021-04-13"
$ws.Range("L7").Value = 'Primera'
$ws.Range("M7").Value = 100
$ws.Range("N7").Value = 15000
$ws.Range("O7").Value = 15000
$ws.Range("P7").Value = 15000
$ws.Range("Q7").Value = '$/caja 15 kilos granel'
$ws.Range("R7").Value = 'Provincia de Curicó'
$ws.Range("T7").Value = 15
$ws.Range("D8").Value = [datetime]"2021-06-07"
$ws.Range("N8").Value = 18000
$ws.Range("O8").Value = 18000
$ws.Range("P8").Value = 18000
$ws.Range("S8").Value = 1000
$ws.Range("D9").Value = [datetime]"2022-05-10"
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 17000
$ws.Range("O9").Value = 17000
$ws.Range("P9").Value = 17000
$ws.Range("S9").Value = 944
$ws.Range("D10").Value = [datetime]"2022-05-11"
$ws.Range("M10").Value = 150
$ws.Range("N10").Value = 17000
$ws.Range("O10").Value = 17000
$ws.Range("P10").Value = 17000
$ws.Range("S10").Value = 944
$ws.Range("D11").Value = [datetime]"2022-04-29"
$ws.Range("R11").Value = 'Provincia de Limarí'
$ws.Range("D12").Value = [datetime]"2021-05-10"
$ws.Range("L12").Value = 'Especial'
$ws.Range("M12").Value = 300
$ws.Range("N12").Value = 20000
$ws.Range("O12").Value = 20000
$ws.Range("P12").Value = 20000
$ws.Range("Q12").Value = '$/caja 18 kilos granel'
$ws.Range("R12").Value = 'Provincia de Limarí'
$ws.Range("S12").Value = 1111
$ws.Range("T12").Value = 18
$ws.Range("D13").Value = [datetime]"2021-06-08"
$ws.Range("M13").Value = 50
$ws.Range("N13").Value = 18000
$ws.Range("O13").Value = 18000
$ws.Range("P13").Value = 18000
$ws.Range("S13").Value = 1000
$ws.Range("D14").Value = [datetime]"2021-04-05"
$ws.Range("L14").Value = 'Primera'
$ws.Range("M14").Value = 150
$ws.Range("N14").Value = 12000
$ws.Range("O14").Value = 12000
$ws.Range("P14").Value = 12000
$ws.Range("Q14").Value = '$/caja 15 kilos granel'
$ws.Range("R14").Value = 'Región Metropolitana'
$ws.Range("S14").Value = 800
$ws.Range("T14").Value = 15
$ws.Range("D15").Value = [datetime]"2023-05-16"
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 200
$ws.Range("N15").Value = 15000
$ws.Range("O15").Value = 15000
$ws.Range("P15").Value = 15000
$ws.Range("Q15").Value = '$/caja 15 kilos granel'
$ws.Range("R15").Value = 'Provincia de Curicó'
$ws.Range("S15").Value = 1000
$ws.Range("T15").Value = 15
$ws.Range("D16").Value = [datetime]"2023-05-08"
$ws.Range("M16").Value = 200
$ws.Range("N16").Value = 14000
$ws.Range("O16").Value = 14000
$ws.Range("P16").Value = 14000
$ws.Range("Q16").Value = '$/caja 18 kilos granel'
$ws.Range("S16").Value = 778
$ws.Range("T16").Value = 18
$ws.Range("D17").Value = [datetime]"2021-06-01"
$ws.Range("L17").Value = 'Especial'
$ws.Range("M17").Value = 200
$ws.Range("N17").Value = 20000
$ws.Range("O17").Value = 20000
$ws.Range("P17").Value = 20000
$ws.Range("Q17").Value = '$/caja 18 kilos granel'
$ws.Range("R17").Value = 'Provincia de Limarí'
$ws.Range("S17").Value = 1111
$ws.Range("T17").Value = 18
$ws.Range("D18").Value = [datetime]"2021-05-26"
$ws.Range("L18").Value = 'Especial'
$ws.Range("M18").Value = 300
$ws.Range("N18").Value = 20000
$ws.Range("O18").Value = 20000
$ws.Range("P18").Value = 20000
$ws.Range("S18").Value = 1111
$ws.Range("D19").Value = [datetime]"2021-05-24"
$ws.Range("M19").Value = 230
$ws.Range("N19").Value = 20000
$ws.Range("O19").Value = 20000
$ws.Range("P19").Value = 20000
$ws.Range("Q19").Value = '$/caja 18 kilos granel'
$ws.Range("S19").Value = 1111
$ws.Range("T19").Value = 18
$ws.Range("D20").Value = [datetime]"2022-06-02"
$ws.Range("M20").Value = 100
$ws.Range("D21").Value = [datetime]"2021-05-03"
$ws.Range("L21").Value = 'Especial'
$ws.Range("M21").Value = 120
